$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pickup/dropoff dates (shared strings) on row 2
$ws.Range("C2").Value = "12/12/2017"
$ws.Range("D2").Value = "12/24/2017"

# Move the active selection to D2 (matches the author's final cursor position)
$ws.Range("D2").Select() | Out-Null
